# Applies the "Finished some of tasks from Excel file and added a few new ones"
# commit: marks several existing tasks as moved/completed, clears their notes,
# and appends five brand-new tasks (rows 14-18) plus three blank spacer rows
# (19-21) to the tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1. Update status / notes for existing tasks 6-10 (rows 9-13)
# ---------------------------------------------------------------------------

# Task 6 - "Email z zaproszeniem" -> moved to task 11
$ws.Range("F9").Value = "Przeniesione"
$ws.Range("G9").Value = "Przeniesione do zadania 11"

# Task 7 - "Tabela z emailami" -> moved to task 11
$ws.Range("F10").Value = "Przeniesione"
$ws.Range("G10").Value = "Przeniesione do zadania 11"

# Task 8 - "Refactoring kontrolerow" -> moved to task 14
$ws.Range("F11").Value = "Przeniesione"
$ws.Range("G11").Value = "Przeniesione do zadania 14"

# Task 9 - "Usuwanie zmienionego hasla" -> finished, clear notes
$ws.Range("F12").Value = "Zakończone"
$ws.Range("G12").ClearContents()

# Task 10 - "Modyfikacja bazy danych" -> finished, clear notes
$ws.Range("F13").Value = "Zakończone"
$ws.Range("G13").ClearContents()

# ---------------------------------------------------------------------------
# 2. Add the new tasks 11-15 (rows 14-18), reusing formatting from existing
#    rows so no new cell styles get created.
# ---------------------------------------------------------------------------

# Row 14 (task 11) - same look as row 9 (C: centered number, D-G: wrapped)
$ws.Range("C9:G9").Copy()
$ws.Range("C14:G14").PasteSpecial($xlPasteFormats)
$ws.Range("C14").Value = 11
$ws.Range("D14").Value = "Uzupełnienie tabeli z emailami"
$ws.Range("E14").Value = "Należy przenieść wszystkie istniejące oraz dodać te brakujące emaile, które obecnie znajdują się w statycznym repozytorium, do tabeli na bazie danych."
$ws.Range("F14").Value = "Nie rozpoczęte"
$ws.Range("G14").Value = "-"
$ws.Rows.Item(14).RowHeight = 72

# Row 15 (task 12)
$ws.Range("C9:G9").Copy()
$ws.Range("C15:G15").PasteSpecial($xlPasteFormats)
$ws.Range("C15").Value = 12
$ws.Range("D15").Value = "Wyczyszczenie repozytoriów"
$ws.Range("E15").Value = "W momencie kiedy nie będą już potrzebne, należy wyczyścić statyczne repozytoria oraz wszystkie odniesienia do nich w kodzie"
$ws.Range("F15").Value = "Nie rozpoczęte"
$ws.Range("G15").Value = "-"
$ws.Rows.Item(15).RowHeight = 57.6

# Row 16 (task 13) - here even column C takes the wrapped style (s=2)
$ws.Range("D9:G9").Copy()
$ws.Range("C16").PasteSpecial($xlPasteFormats)
$ws.Range("D16:G16").PasteSpecial($xlPasteFormats)
$ws.Range("C16").Value = 13
$ws.Range("D16").Value = "Usuwanie filmów"
$ws.Range("E16").Value = "Dodać metodę do usuwania filmów z bazy"
$ws.Range("F16").Value = "Nie rozpoczęte"
$ws.Range("G16").Value = "-"
$ws.Rows.Item(16).RowHeight = 28.8

# Row 17 (task 14)
$ws.Range("D9:G9").Copy()
$ws.Range("C17").PasteSpecial($xlPasteFormats)
$ws.Range("D17:G17").PasteSpecial($xlPasteFormats)
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = "Refactoring Music Controllera"
$ws.Range("E17").Value = "Należy usunąć wszystkie nieużywane i przestarzałe metody z kontrollera"
$ws.Range("F17").Value = "Nie rozpoczęte"
$ws.Range("G17").Value = "-"
$ws.Rows.Item(17).RowHeight = 28.8

# Row 18 (task 15)
$ws.Range("D9:G9").Copy()
$ws.Range("C18").PasteSpecial($xlPasteFormats)
$ws.Range("D18:G18").PasteSpecial($xlPasteFormats)
$ws.Range("C18").Value = 15
$ws.Range("D18").Value = "Usuwanie muzyki"
$ws.Range("E18").Value = "Dodać metodę do usuwania muzyki z bazy"
$ws.Range("F18").Value = "Nie rozpoczęte"
$ws.Range("G18").Value = "-"
$ws.Rows.Item(18).RowHeight = 28.8

# ---------------------------------------------------------------------------
# 3. Three trailing blank (but formatted) rows 19-21
# ---------------------------------------------------------------------------
$ws.Range("D9:G9").Copy()
$ws.Range("C19").PasteSpecial($xlPasteFormats)
$ws.Range("D19:G19").PasteSpecial($xlPasteFormats)

$ws.Range("D9:G9").Copy()
$ws.Range("C20").PasteSpecial($xlPasteFormats)
$ws.Range("D20:G20").PasteSpecial($xlPasteFormats)

$ws.Range("D9:G9").Copy()
$ws.Range("C21").PasteSpecial($xlPasteFormats)
$ws.Range("D21:G21").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------------
# 4. Refresh the visible selection/scroll position to match the saved file
# ---------------------------------------------------------------------------
[void]$ws.Range("G19").Select()
